$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 15554.84383846409
$ws.Range("T2").Value = 67
$ws.Range("U2").Value = 130672.5018699766
$ws.Range("V2").Value = 776666.1140470143
$ws.Range("W2").Value = 18882.39501835127
$ws.Range("X2").Value = 438926.6492088101
$ws.Range("Y2").Value = 696517.4566331406
$ws.Range("Z2").Value = 9571.926512375001
$ws.Range("AA2").Value = 456991.2390601542
$ws.Range("AB2").Value = 267533.8774694387
$ws.Range("AC2").Value = "High Growth"
$ws.Range("AI2").Value = 174908.8743095075
$ws.Range("AJ2").Value = 0.09259616511851354
$ws.Range("AK2").Value = 0.009240073486540767
$ws.Range("AL2").Value = 227.5844075966687

# Row 3
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 61533.85532669698
$ws.Range("T3").Value = 67
$ws.Range("U3").Value = 52799.45170986404
$ws.Range("V3").Value = 946917.3434636819
$ws.Range("W3").Value = 12165.05731708026
$ws.Range("X3").Value = 387879.9105327942
$ws.Range("Y3").Value = 131030.7875657655
$ws.Range("Z3").Value = 13541.63148054171
$ws.Range("AA3").Value = 340603.0693258547
$ws.Range("AB3").Value = 404932.7423894164
$ws.Range("AC3").Value = "Balanced"
$ws.Range("AI3").Value = 704474.258393494
$ws.Range("AJ3").Value = 0.2174504927055718
$ws.Range("AK3").Value = 0.8330380975093908
$ws.Range("AL3").Value = 439.4348708887354

# Row 4
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 36799.48332397205
$ws.Range("T4").Value = 67
$ws.Range("U4").Value = 136419.9448274699
$ws.Range("V4").Value = 762499.2863315954
$ws.Range("W4").Value = 11933.10812590106
$ws.Range("X4").Value = 484020.3955052439
$ws.Range("Y4").Value = 103597.4160782679
$ws.Range("Z4").Value = 12140.90122232903
$ws.Range("AA4").Value = 196107.3374299094
$ws.Range("AB4").Value = 153248.8255119704
$ws.Range("AC4").Value = "Balanced"
$ws.Range("AI4").Value = 463150.1999901332
$ws.Range("AJ4").Value = 0.7492544170632397
$ws.Range("AK4").Value = 0.984329434348122
$ws.Range("AL4").Value = 246.1341207525793

# Row 5
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 30789.47658762065
$ws.Range("T5").Value = 67
$ws.Range("U5").Value = 69651.49919650386
$ws.Range("V5").Value = 372433.6284770427
$ws.Range("W5").Value = 15672.88490603777
$ws.Range("X5").Value = 299409.0957879981
$ws.Range("Y5").Value = 513445.819141016
$ws.Range("Z5").Value = 14293.93872474714
$ws.Range("AA5").Value = 283870.0371462243
$ws.Range("AB5").Value = 446133.921401752
$ws.Range("AC5").Value = "High Growth"
$ws.Range("AI5").Value = 840428.5325354325
$ws.Range("AJ5").Value = 0.7316937564646905
$ws.Range("AK5").Value = 0.7034947841347021
$ws.Range("AL5").Value = 357.7805787295199

# Row 6
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 6183.518460005614
$ws.Range("S6").Value = 60015.28804864444
$ws.Range("T6").Value = 67
$ws.Range("U6").Value = 57465.8150156651
$ws.Range("V6").Value = 515071.490046539
$ws.Range("W6").Value = 10000.52595458603
$ws.Range("X6").Value = 145575.7245006389
$ws.Range("Y6").Value = 886954.3340915912
$ws.Range("Z6").Value = 9389.94366182342
$ws.Range("AA6").Value = 50024.91058925386
$ws.Range("AB6").Value = 191161.5766286134
$ws.Range("AC6").Value = "Balanced"
$ws.Range("AI6").Value = 204865.7621804457
$ws.Range("AJ6").Value = 0.04561463715088643
$ws.Range("AK6").Value = 0.1816312004920487
$ws.Range("AL6").Value = 243.0797801555988

# Row 7
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 20791.75864178867
$ws.Range("T7").Value = 67
$ws.Range("U7").Value = 58401.50739581222
$ws.Range("V7").Value = 372035.8609407893
$ws.Range("W7").Value = 1007.401135150876
$ws.Range("X7").Value = 150099.4509740684
$ws.Range("Y7").Value = 977584.1672633904
$ws.Range("Z7").Value = 9120.292464357344
$ws.Range("AA7").Value = 190337.0993069851
$ws.Range("AB7").Value = 480852.8850465292
$ws.Range("AC7").Value = "Balanced"
$ws.Range("AI7").Value = 164958.8680867573
$ws.Range("AJ7").Value = 0.2091570294808411
$ws.Range("AK7").Value = 0.5123934632959037
$ws.Range("AL7").Value = 354.2740750100015

# Row 8
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 63750.81939309922
$ws.Range("T8").Value = 62
$ws.Range("U8").Value = 51011.76230108475
$ws.Range("V8").Value = 220515.7478719629
$ws.Range("W8").Value = 13981.96149437022
$ws.Range("X8").Value = 148437.2181815473
$ws.Range("Y8").Value = 636775.0202910862
$ws.Range("Z8").Value = 18128.36174120342
$ws.Range("AA8").Value = 241453.1920784792
$ws.Range("AB8").Value = 262288.2069017605
$ws.Range("AC8").Value = "High Growth"
$ws.Range("AI8").Value = 124833.057490086
$ws.Range("AJ8").Value = 0.2869150433006046
$ws.Range("AK8").Value = 0.5804471371881301
$ws.Range("AL8").Value = 249.070112753735

# Row 9
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 677.8165367962301
$ws.Range("S9").Value = 55139.94585884939
$ws.Range("T9").Value = 67
$ws.Range("U9").Value = 57309.9221252682
$ws.Range("V9").Value = 212437.2687489749
$ws.Range("W9").Value = 19847.92797777265
$ws.Range("X9").Value = 306308.0905586321
$ws.Range("Y9").Value = 565373.8362598775
$ws.Range("Z9").Value = 2744.408401943892
$ws.Range("AA9").Value = 448401.9468243135
$ws.Range("AB9").Value = 370212.7176773662
$ws.Range("AC9").Value = "Growth"
$ws.Range("AI9").Value = 722080.6624027486
$ws.Range("AJ9").Value = 0.6772633053384269
$ws.Range("AK9").Value = 0.7875421559550259
$ws.Range("AL9").Value = 422.2750397359712

# Row 10
$ws.Range("Q10").Value = 9222.339545700042
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 2375.362006503465
$ws.Range("T10").Value = 67
$ws.Range("U10").Value = 50332.55947467734
$ws.Range("V10").Value = 490278.4285152588
$ws.Range("W10").Value = 5345.250757629565
$ws.Range("X10").Value = 253449.0657650416
$ws.Range("Y10").Value = 861882.9012951787
$ws.Range("Z10").Value = 4584.386461731588
$ws.Range("AA10").Value = 355945.7553830463
$ws.Range("AB10").Value = 119162.4372535084
$ws.Range("AC10").Value = "Growth"
$ws.Range("AI10").Value = 30452.96494160032
$ws.Range("AJ10").Value = 0.06303828966937675
$ws.Range("AK10").Value = 0.6064754600685757
$ws.Range("AL10").Value = 97.1711114714916

# Row 11
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 13381.0203582225
$ws.Range("T11").Value = 67
$ws.Range("U11").Value = 138958.5234374146
$ws.Range("V11").Value = 880844.3097627467
$ws.Range("W11").Value = 13581.81231254828
$ws.Range("X11").Value = 486606.5073932452
$ws.Range("Y11").Value = 756642.8180442983
$ws.Range("Z11").Value = 17631.70798421161
$ws.Range("AA11").Value = 255258.3974552879
$ws.Range("AB11").Value = 378698.9794550575
$ws.Range("AC11").Value = "Balanced"
$ws.Range("AI11").Value = 746994.2401995035
$ws.Range("AJ11").Value = 0.5556492427701452
$ws.Range("AK11").Value = 0.2184028548119228
$ws.Range("AL11").Value = 386.6630048052849
